# "Analises de linha do tempo"
# - A10:A13 switch from inline-string labels ("12".."15") to the correct
#   sequential numeric index (8..11), keeping their existing style.
# - Four new rows (14-17) are appended, mirroring the existing 4-city
#   weather block (Fortaleza / Sobral / Acaraú / Itarema) for 12\09\24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mislabeled index column for the existing last block ---------
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(13, 1).Value = 11

# --- Bring formatting (border/bold/center style) for the new rows --------
$ws.Range("A10:M13").Copy()
$ws.Range("A14:M17").PasteSpecial(-4122)

function Set-TextCellA($row, $text) {
    # Column A holds small integer-looking labels ("12".."15") that Excel's
    # auto-detection would otherwise coerce to numbers, so force text,
    # write the value, then re-apply the bordered/bold/centered style that
    # the format-paste already put on this cell (NumberFormat overwrote it).
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-TextCell($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = $text
}

function Set-NumCell($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# --- Row 14: Fortaleza -----------------------------------------------------
Set-TextCellA 14 "12"
Set-TextCell 14 2 "12\09\24"
Set-TextCell 14 3 "qui"
Set-TextCell 14 4 "13:32:57"
Set-TextCell 14 5 "Fortaleza"
Set-TextCell 14 6 "algumas nuvens"
Set-NumCell  14 7 31.15
Set-NumCell  14 8 31.07
Set-NumCell  14 9 31.15
Set-NumCell  14 10 38.15
Set-NumCell  14 11 82
Set-NumCell  14 12 1013
Set-NumCell  14 13 10.29

# --- Row 15: Sobral ----------------------------------------------------
Set-TextCellA 15 "13"
Set-TextCell 15 2 "12\09\24"
Set-TextCell 15 3 "qui"
Set-TextCell 15 4 "13:32:57"
Set-TextCell 15 5 "Sobral"
Set-TextCell 15 6 "céu limpo"
Set-NumCell  15 7 36.64
Set-NumCell  15 8 36.64
Set-NumCell  15 9 36.64
Set-NumCell  15 10 35.93
Set-NumCell  15 11 25
Set-NumCell  15 12 1010
Set-NumCell  15 13 4.16

# --- Row 16: Acaraú ----------------------------------------------------
Set-TextCellA 16 "14"
Set-TextCell 16 2 "12\09\24"
Set-TextCell 16 3 "qui"
Set-TextCell 16 4 "13:32:57"
Set-TextCell 16 5 "Acaraú"
Set-TextCell 16 6 "céu limpo"
Set-NumCell  16 7 29.81
Set-NumCell  16 8 29.81
Set-NumCell  16 9 29.81
Set-NumCell  16 10 32.1
Set-NumCell  16 11 58
Set-NumCell  16 12 1011
Set-NumCell  16 13 10.51

# --- Row 17: Itarema ---------------------------------------------------
Set-TextCellA 17 "15"
Set-TextCell 17 2 "12\09\24"
Set-TextCell 17 3 "qui"
Set-TextCell 17 4 "13:32:57"
Set-TextCell 17 5 "Itarema"
Set-TextCell 17 6 "céu limpo"
Set-NumCell  17 7 29.48
Set-NumCell  17 8 29.48
Set-NumCell  17 9 29.48
Set-NumCell  17 10 31.53
Set-NumCell  17 11 58
Set-NumCell  17 12 1012
Set-NumCell  17 13 9.77

# --- Restore the bordered/bold/center style on column A (NumberFormat
#     above blew it away) by pasting just the formats from A13 again.
$ws.Range("A13").Copy()
$ws.Range("A14:A17").PasteSpecial(-4122)

Write-Output "Applied: A10:A13 renumbered, rows 14-17 appended."
